# Update MSME summary figures for Egypt, Arab Rep. with more precise values.
# The source cells are stored as text (not numbers), so a leading apostrophe
# is used to force text entry; the cell Style is then re-applied from an
# already-clean neighboring text cell so the quote-prefix formatting that
# Excel would otherwise stamp on the cell is not left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $value
    $range.Style = $ws.Range("B12").Style
}

# Enterprises density (per 1000 people)
Set-TextValue "B13" "33.53"
Set-TextValue "C13" "0.08"
Set-TextValue "D13" "33.61"

# Employment (% of total)
Set-TextValue "B14" "74.29"
Set-TextValue "C14" "20.53"
Set-TextValue "D14" "94.82"

# Enterprises (% of total)
Set-TextValue "B16" "99.73"
Set-TextValue "C16" "0.25"
Set-TextValue "D16" "99.98"
